# Correção das notas do fórum para MATC65 em 2021.2
# Zera todas as notas semanais e a nota_forum final (colunas B:Q, linhas 2:50)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:Q50").Value = 0
